$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "Booleans" column (F) — header + two boolean cells for rows 2/3,
# matching the "add boolvalue to parser" fixture update.
$ws.Range("F1").Value = "Booleans"

$ws.Range("F2").Value = $true
$ws.Range("F3").Value = $false

# Boolean-style display format ("TRUE"/"FALSE") applied to the new cells.
$ws.Range("F2:F3").NumberFormat = '"TRUE";"TRUE";"FALSE"'

$ws.Range("F4").Select()
